$d = $word.ActiveDocument

# --- 1. Replace the text of the "Data" bullet (paragraph 2) ---
$oldData = "Data: Attributes / Values x Types / Instances scaling. Price / Amount, Product / Item. Measures: (Attribute, Value)."
$newData = "Data: Attributes / Values x Types / Instances scaling. Price / Amount, Product / Item. Measures: (Attributes, Values)."
$d.Content.Find.Execute($oldData, $true, $false, $false, $false, $false, $true, 1, $false, $newData, 2) | Out-Null

# --- 2. Replace the text of the "Information" bullet (paragraph 3) ---
$oldInfo = "Information: Data across Dimensions Attributes / Values scaling. Time / Date, POS / Store Attributes Variation. Dimensions: (Dimension, Value)."
$newInfo = "Information: Data across Dimensions Attributes / Values scaling. Time / Date / Price, POS / Store / Availability Attributes (Variation). Dimensions: (Dimension Type, Dimension Value)."
$d.Content.Find.Execute($oldInfo, $true, $false, $false, $false, $false, $true, 1, $false, $newInfo, 2) | Out-Null

# --- 3. Promote the "Data", "Information" and "Knowledge" bullets one list level
#        (ilvl 0 -> 1) and re-indent them (left 600 -> 1440 twips = 30pt -> 72pt,
#        hanging stays 360 twips = 18pt). ---
for ($i = 2; $i -le 4; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.ListFormat.ListLevelNumber = 2
    $p.Range.ParagraphFormat.LeftIndent = 72
    $p.Range.ParagraphFormat.FirstLineIndent = -18
}

# --- 4. Insert a new, empty, non-list paragraph right after the "Knowledge"
#        bullet (paragraph 4). ---
$knowledge = $d.Paragraphs.Item(4)
$knowledge.Range.InsertParagraphAfter()
$blank = $d.Paragraphs.Item(5)
$blank.Range.ListFormat.RemoveNumbers()
$blank.Range.ParagraphFormat.LeftIndent = 0
$blank.Range.ParagraphFormat.FirstLineIndent = 0
